# Add "ROIs matching" step to the parcellation stage.
#
# Inserts two new parameter rows after "parcellation.templateScript"
# (current row 31):
#   - parcellation.matchROIs
#   - parcellation.lutFile
# and removes the now-redundant "collect_region_properties.lutFile" row,
# since the lookup-table file parameter moved to the parcellation step.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 32 (pushes everything from the old row 32
# down to row 34).
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(32).Insert()

# New row 32: parcellation.matchROIs
$ws.Range("A32").Value = "parcellation.matchROIs"
$ws.Range("D32").Value = "parcellation"
$ws.Range("E32").Value = "logical"
$ws.Range("G32").Value = "standard"
$ws.Range("H32").Value = "Flag whether the parcellation step should reassign the ROIs in the parcellationFile to match the template's color lookup table."

# New row 33: parcellation.lutFile
$ws.Range("A33").Value = "parcellation.lutFile"
$ws.Range("B33").Value = "parcellation,collect_region_properties"
$ws.Range("E33").Value = "char "
$ws.Range("F33").Value = "isfile nonempty"
$ws.Range("G33").Value = "standard"
$ws.Range("H33").Value = "Freesurfer's color lookup table of the template"

# Remove the now-redundant collect_region_properties.lutFile row. It used
# to be row 35, but the two inserted rows above shifted it down to row 37.
$ws.Rows.Item(37).Delete()

# Restore view state: selection + scroll position.
$ws.Range("A8").Select()
$ws.Range("B33").Select()
